$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Coin/Link/Price/Volume columns so numeric-looking
# strings (e.g. "224.90", "1.00") are preserved exactly as text, matching
# the original inlineStr cell contents instead of being parsed as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '36.373.88'
$ws.Cells.Item(2, 5).Value = '  -3.13%  '
$ws.Cells.Item(3, 4).Value = '1.946.60'
$ws.Cells.Item(3, 5).Value = '  -4.00%  '
$ws.Cells.Item(4, 5).Value = '  +0.14%  '
$ws.Cells.Item(5, 4).Value = '224.90'
$ws.Cells.Item(5, 5).Value = '  -12.05%  '
$ws.Cells.Item(6, 4).Value = '0.587'
$ws.Cells.Item(6, 5).Value = '  -4.50%  '
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 4).Value = '52.16'
$ws.Cells.Item(8, 5).Value = '  -8.67%  '
$ws.Cells.Item(9, 4).Value = '0.359'
$ws.Cells.Item(9, 5).Value = '  -6.43%  '
$ws.Cells.Item(10, 4).Value = '56.55'
$ws.Cells.Item(10, 5).Value = '  -1.21%  '
$ws.Cells.Item(11, 4).Value = '0.0722'
$ws.Cells.Item(11, 5).Value = '  -8.14%  '
$ws.Cells.Item(12, 4).Value = '0.0963'
$ws.Cells.Item(12, 5).Value = '  -5.14%  '
$ws.Cells.Item(13, 4).Value = '2.231.34'
$ws.Cells.Item(13, 5).Value = '  -3.88%  '
$ws.Cells.Item(14, 4).Value = '13.50'
$ws.Cells.Item(14, 5).Value = '  -7.04%  '
$ws.Cells.Item(15, 4).Value = '19.33'
$ws.Cells.Item(15, 5).Value = '  -8.83%  '
$ws.Cells.Item(16, 4).Value = '0.730'
$ws.Cells.Item(16, 5).Value = '  -10.73%  '
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '1.958.46'
$ws.Cells.Item(17, 5).Value = '  -3.76%  '
$ws.Cells.Item(18, 2).Value = 'Polkadot'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(18, 4).Value = '4.90'
$ws.Cells.Item(18, 5).Value = '  -8.62%  '
$ws.Cells.Item(19, 4).Value = '36.248.27'
$ws.Cells.Item(19, 5).Value = '  -3.21%  '
$ws.Cells.Item(20, 4).Value = '66.02'
$ws.Cells.Item(20, 5).Value = '  -5.12%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0777'
$ws.Cells.Item(21, 5).Value = '  -8.69%  '
$ws.Cells.Item(22, 4).Value = '4.89'
$ws.Cells.Item(22, 5).Value = '  -6.15%  '
$ws.Cells.Item(23, 4).Value = '217.27'
$ws.Cells.Item(23, 5).Value = '  -5.06%  '
$ws.Cells.Item(24, 5).Value = '  -0.10%  '
$ws.Cells.Item(25, 4).Value = '2.31'
$ws.Cells.Item(25, 5).Value = '  -1.15%  '
$ws.Cells.Item(26, 4).Value = '2.26'
$ws.Cells.Item(26, 5).Value = '  -14.19%  '
$ws.Cells.Item(27, 4).Value = '158.88'
$ws.Cells.Item(27, 5).Value = '  -3.19%  '
$ws.Cells.Item(28, 4).Value = '8.34'
$ws.Cells.Item(28, 5).Value = '  -7.95%  '
$ws.Cells.Item(29, 4).Value = '18.39'
$ws.Cells.Item(29, 5).Value = '  -7.63%  '
$ws.Cells.Item(30, 2).Value = 'ImmutableX'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(30, 4).Value = '1.27'
$ws.Cells.Item(30, 5).Value = '  -8.13%  '
$ws.Cells.Item(31, 2).Value = 'Kaspa'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(31, 4).Value = '0.117'
$ws.Cells.Item(31, 5).Value = '  -11.04%  '
$ws.Cells.Item(32, 4).Value = '0.114'
$ws.Cells.Item(32, 5).Value = '  -5.28%  '
$ws.Cells.Item(33, 5).Value = '  -10.45%  '
$ws.Cells.Item(34, 4).Value = '0.0590'
$ws.Cells.Item(34, 5).Value = '  -11.56%  '
$ws.Cells.Item(35, 4).Value = '4.09'
$ws.Cells.Item(35, 5).Value = '  -10.61%  '
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).Value = '2.22'
$ws.Cells.Item(36, 5).Value = '  -8.82%  '
$ws.Cells.Item(37, 2).Value = 'BinanceUSD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(37, 4).Value = '1.00'
$ws.Cells.Item(37, 5).Value = '  -0.07%  '
$ws.Cells.Item(38, 5).Value = '  -2.19%  '
$ws.Cells.Item(39, 4).Value = '3.11'
$ws.Cells.Item(39, 5).Value = '  -8.76%  '
$ws.Cells.Item(40, 4).Value = '3.05'
$ws.Cells.Item(40, 5).Value = '  +0.54%  '
$ws.Cells.Item(41, 4).Value = '4.97'
$ws.Cells.Item(41, 5).Value = '  -7.23%  '
$ws.Cells.Item(42, 4).Value = '1.387.52'
$ws.Cells.Item(42, 5).Value = '  -1.42%  '
$ws.Cells.Item(43, 4).Value = '0.0196'
$ws.Cells.Item(43, 5).Value = '  -9.17%  '
$ws.Cells.Item(44, 4).Value = '0.0854'
$ws.Cells.Item(44, 5).Value = '  -11.77%  '
$ws.Cells.Item(45, 4).Value = '1.05'
$ws.Cells.Item(45, 5).Value = '  -12.39%  '
$ws.Cells.Item(46, 4).Value = '84.63'
$ws.Cells.Item(46, 5).Value = '  -7.01%  '
$ws.Cells.Item(47, 2).Value = 'MXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(47, 4).Value = '2.84'
$ws.Cells.Item(47, 5).Value = '  -1.18%  '
$ws.Cells.Item(48, 2).Value = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).Value = '0.961'
$ws.Cells.Item(48, 5).Value = '  -7.49%  '
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '14.39'
$ws.Cells.Item(49, 5).Value = '  -10.31%  '
$ws.Cells.Item(50, 5).Value = '  -9.81%  '
$ws.Cells.Item(51, 4).Value = '2.125.06'
$ws.Cells.Item(51, 5).Value = '  -3.94%  '
